$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 61; $r -le 582; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value2 = $cell.Value2 * 1000
}
